$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 188 ---
# Column A (date) needs the same date-time number format / style as the
# rows above it (s="1" -> numFmtId 165 "yyyy-mm-dd hh:mm:ss"). Copy the
# format from the last existing data row (A187) instead of re-creating a
# new style entry.
$ws.Range("A187").Copy()
$ws.Range("A188").PasteSpecial(-4122)
$ws.Range("A188").Value = 45468.2916666667

$ws.Range("B188").Value = 0
$ws.Range("C188").Value = 2.83999991416931
$ws.Range("D188").Value = 2.83999991416931
$ws.Range("E188").Value = 2.83999991416931
$ws.Range("F188").Value = 2.83999991416931

# Column G ("adj_close") stores the close price formatted as text so it
# lands in the shared-string table (matches the existing rows). Use a
# leading apostrophe to force text, then reset the cell style back to
# Normal so no stray number-format style sticks around on the cell.
$ws.Range("G188").Value = "'2.83999991416931"
$ws.Range("G188").Style = "Normal"

$ws.Range("H188").Value = "'EAV.MI"
$ws.Range("H188").Style = "Normal"

# --- Row 189 ---
$ws.Range("A187").Copy()
$ws.Range("A189").PasteSpecial(-4122)
$ws.Range("A189").Value = 45469.6423726852

$ws.Range("B189").Value = 3500
$ws.Range("C189").Value = 2.85999989509583
$ws.Range("D189").Value = 2.76999998092651
$ws.Range("E189").Value = 2.85999989509583
$ws.Range("F189").Value = 2.77999997138977

$ws.Range("G189").Value = "'2.77999997138977"
$ws.Range("G189").Style = "Normal"

$ws.Range("H189").Value = "'EAV.MI"
$ws.Range("H189").Style = "Normal"

$excel.CutCopyMode = 0
